$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(12, 8).Value = 273.88235
$ws.Cells.Item(12, 9).Value = 243.66667
$ws.Cells.Item(12, 10).Value = 500.5
$ws.Cells.Item(12, 11).Value = 243.66667
$ws.Cells.Item(12, 12).Value = 500.5
$ws.Cells.Item(12, 13).Value = -73.66667000000001
$ws.Cells.Item(12, 14).Value = -840.5
$ws.Cells.Item(18, 8).Value = 1477.1666
$ws.Cells.Item(18, 9).Value = 1234.6471
$ws.Cells.Item(18, 10).Value = 5600
$ws.Cells.Item(18, 11).Value = 1234.6471
$ws.Cells.Item(18, 12).Value = 5600
$ws.Cells.Item(18, 13).Value = -950.6470999999999
$ws.Cells.Item(18, 14).Value = -6168
$ws.Cells.Item(40, 8).Value = 3277.1875
$ws.Cells.Item(40, 9).Value = 2680.3333
$ws.Cells.Item(40, 10).Value = 4044.5715
$ws.Cells.Item(40, 11).Value = 2680.3333
$ws.Cells.Item(40, 12).Value = 4044.5715
$ws.Cells.Item(40, 13).Value = -2505.3333
$ws.Cells.Item(40, 14).Value = -4394.5715
$ws.Cells.Item(98, 8).Value = 62501304
$ws.Cells.Item(98, 9).Value = 71429670
$ws.Cells.Item(98, 11).Value = 71429670
$ws.Cells.Item(98, 13).Value = -71428172
$ws.Cells.Item(107, 8).Value = 347.27274
$ws.Cells.Item(107, 9).Value = 372.3
$ws.Cells.Item(107, 11).Value = 372.3
$ws.Cells.Item(107, 13).Value = 1547.7
$ws.Cells.Item(122, 8).Value = 62501304
$ws.Cells.Item(122, 9).Value = 71429670
$ws.Cells.Item(122, 11).Value = 214289010
$ws.Cells.Item(122, 13).Value = -214286560
$ws.Cells.Item(127, 8).Value = 2274.375
$ws.Cells.Item(127, 9).Value = 456.42856
$ws.Cells.Item(127, 11).Value = 1369.28568
$ws.Cells.Item(127, 13).Value = 3590.71432
$ws.Cells.Item(138, 8).Value = 2503.2627
$ws.Cells.Item(138, 9).Value = 1608
$ws.Cells.Item(138, 10).Value = 2615.1704
$ws.Cells.Item(138, 11).Value = 4824
$ws.Cells.Item(138, 12).Value = 7845.5112
$ws.Cells.Item(138, 13).Value = 316
$ws.Cells.Item(138, 14).Value = -18125.5112

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(4, 8).Value = 533.8333
$ws.Cells.Item(4, 9).Value = 575.25
$ws.Cells.Item(4, 11).Value = 575.25
$ws.Cells.Item(4, 13).Value = -459.25
$ws.Cells.Item(5, 8).Value = 239.46666
$ws.Cells.Item(5, 9).Value = 98.72727
$ws.Cells.Item(5, 10).Value = 626.5
$ws.Cells.Item(5, 11).Value = 98.72727
$ws.Cells.Item(5, 12).Value = 626.5
$ws.Cells.Item(5, 13).Value = 13.27273
$ws.Cells.Item(5, 14).Value = -850.5
$ws.Cells.Item(21, 8).Value = 5250
$ws.Cells.Item(21, 10).Value = 5500
$ws.Cells.Item(21, 12).Value = 5500
$ws.Cells.Item(21, 14).Value = -6248
$ws.Cells.Item(61, 8).Value = 50010476
$ws.Cells.Item(61, 9).Value = 38471916
$ws.Cells.Item(61, 11).Value = 38471916
$ws.Cells.Item(61, 13).Value = -38471704
$ws.Cells.Item(74, 8).Value = 17340678
$ws.Cells.Item(74, 9).Value = 50006400
$ws.Cells.Item(74, 11).Value = 50006400
$ws.Cells.Item(74, 13).Value = -50005526
$ws.Cells.Item(77, 8).Value = 17340678
$ws.Cells.Item(77, 9).Value = 50006400
$ws.Cells.Item(77, 11).Value = 250032000
$ws.Cells.Item(77, 13).Value = -250027632
$ws.Cells.Item(96, 8).Value = 28399.8
$ws.Cells.Item(96, 10).Value = 28399.8
$ws.Cells.Item(96, 12).Value = 28399.8
$ws.Cells.Item(96, 14).Value = -33891.8
$ws.Cells.Item(97, 8).Value = 2694
$ws.Cells.Item(97, 9).Value = 2946.3333
$ws.Cells.Item(97, 11).Value = 2946.3333
$ws.Cells.Item(97, 13).Value = -2450.3333
$ws.Cells.Item(118, 8).Value = 52500
$ws.Cells.Item(118, 10).Value = 52500
$ws.Cells.Item(118, 12).Value = 52500
$ws.Cells.Item(118, 14).Value = -55814
$ws.Cells.Item(132, 8).Value = 5874.171
$ws.Cells.Item(132, 9).Value = 3739.5483
$ws.Cells.Item(132, 11).Value = 11218.6449
$ws.Cells.Item(132, 13).Value = -8688.644899999999
$ws.Cells.Item(136, 8).Value = 50010476
$ws.Cells.Item(136, 9).Value = 38471916
$ws.Cells.Item(136, 11).Value = 115415748
$ws.Cells.Item(136, 13).Value = -115413198

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(4, 8).Value = 239.46666
$ws.Cells.Item(4, 9).Value = 98.72727
$ws.Cells.Item(4, 10).Value = 626.5
$ws.Cells.Item(4, 11).Value = 98.72727
$ws.Cells.Item(4, 12).Value = 626.5
$ws.Cells.Item(4, 13).Value = 16.27273
$ws.Cells.Item(4, 14).Value = -856.5
$ws.Cells.Item(19, 8).Value = 6776.923
$ws.Cells.Item(19, 9).Value = 675
$ws.Cells.Item(19, 11).Value = 675
$ws.Cells.Item(19, 13).Value = -502
$ws.Cells.Item(94, 8).Value = 2611
$ws.Cells.Item(94, 9).Value = 9000
$ws.Cells.Item(94, 10).Value = 1812.375
$ws.Cells.Item(94, 11).Value = 9000
$ws.Cells.Item(94, 12).Value = 1812.375
$ws.Cells.Item(94, 13).Value = -8549
$ws.Cells.Item(94, 14).Value = -2714.375
$ws.Cells.Item(107, 8).Value = 1002.25
$ws.Cells.Item(107, 9).Value = 1018.6875
$ws.Cells.Item(107, 11).Value = 1018.6875
$ws.Cells.Item(107, 13).Value = 901.3125

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(17, 8).Value = 0
$ws.Cells.Item(17, 9).Value = 0
$ws.Cells.Item(17, 11).Value = 0
$ws.Cells.Item(17, 13).ClearContents()
$ws.Cells.Item(22, 8).Value = 445.9
$ws.Cells.Item(22, 9).Value = 373.16666
$ws.Cells.Item(22, 10).Value = 555
$ws.Cells.Item(22, 11).Value = 373.16666
$ws.Cells.Item(22, 12).Value = 555
$ws.Cells.Item(22, 13).Value = -23.16665999999998
$ws.Cells.Item(22, 14).Value = -1255
$ws.Cells.Item(25, 8).Value = 0
$ws.Cells.Item(25, 10).Value = 0
$ws.Cells.Item(25, 12).Value = 0
$ws.Cells.Item(25, 14).ClearContents()
$ws.Cells.Item(99, 8).Value = 3376.875
$ws.Cells.Item(99, 9).Value = 4171.6665
$ws.Cells.Item(99, 10).Value = 992.5
$ws.Cells.Item(99, 11).Value = 4171.6665
$ws.Cells.Item(99, 12).Value = 992.5
$ws.Cells.Item(99, 13).Value = -2673.6665
$ws.Cells.Item(99, 14).Value = -3988.5
$ws.Cells.Item(126, 8).Value = 3376.875
$ws.Cells.Item(126, 9).Value = 4171.6665
$ws.Cells.Item(126, 10).Value = 992.5
$ws.Cells.Item(126, 11).Value = 12514.9995
$ws.Cells.Item(126, 12).Value = 2977.5
$ws.Cells.Item(126, 13).Value = -10044.9995
$ws.Cells.Item(126, 14).Value = -7917.5

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(17, 8).Value = 529.2857
$ws.Cells.Item(17, 9).Value = 52.5
$ws.Cells.Item(17, 10).Value = 720
$ws.Cells.Item(17, 11).Value = 157.5
$ws.Cells.Item(17, 12).Value = 2160
$ws.Cells.Item(17, 13).Value = 11.5
$ws.Cells.Item(17, 14).Value = -2498
$ws.Cells.Item(26, 8).Value = 174
$ws.Cells.Item(26, 9).Value = 174
$ws.Cells.Item(26, 11).Value = 522
$ws.Cells.Item(26, 13).Value = -234
$ws.Cells.Item(87, 8).Value = 5432.6665
$ws.Cells.Item(87, 9).Value = 5432.6665
$ws.Cells.Item(87, 11).Value = 16297.9995
$ws.Cells.Item(87, 13).Value = -15049.9995
$ws.Cells.Item(90, 8).Value = 5432.6665
$ws.Cells.Item(90, 9).Value = 5432.6665
$ws.Cells.Item(90, 11).Value = 48893.9985
$ws.Cells.Item(90, 13).Value = -42653.9985
$ws.Cells.Item(118, 8).Value = 7349.75
$ws.Cells.Item(118, 9).Value = 5999
$ws.Cells.Item(118, 10).Value = 7800
$ws.Cells.Item(118, 11).Value = 17997
$ws.Cells.Item(118, 12).Value = 23400
$ws.Cells.Item(118, 13).Value = -16754
$ws.Cells.Item(118, 14).Value = -25886
$ws.Cells.Item(137, 8).Value = 5337.55
$ws.Cells.Item(137, 9).Value = 5599.5
$ws.Cells.Item(137, 10).Value = 5272.0625
$ws.Cells.Item(137, 11).Value = 16798.5
$ws.Cells.Item(137, 12).Value = 15816.1875
$ws.Cells.Item(137, 13).Value = -11698.5
$ws.Cells.Item(137, 14).Value = -26016.1875

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(2, 8).Value = 156.65
$ws.Cells.Item(2, 9).Value = 45.7
$ws.Cells.Item(2, 10).Value = 267.6
$ws.Cells.Item(2, 11).Value = 45.7
$ws.Cells.Item(2, 12).Value = 267.6
$ws.Cells.Item(2, 13).Value = 67.3
$ws.Cells.Item(2, 14).Value = -493.6
$ws.Cells.Item(11, 8).Value = 4079336.8
$ws.Cells.Item(11, 9).Value = 4628000
$ws.Cells.Item(11, 10).Value = 3622117.2
$ws.Cells.Item(11, 11).Value = 4628000
$ws.Cells.Item(11, 12).Value = 3622117.2
$ws.Cells.Item(11, 13).Value = -4627861
$ws.Cells.Item(11, 14).Value = -3622395.2
$ws.Cells.Item(102, 8).Value = 2768.4092
$ws.Cells.Item(102, 9).Value = 1861.3334
$ws.Cells.Item(102, 11).Value = 1861.3334
$ws.Cells.Item(102, 13).Value = -239.3334
$ws.Cells.Item(122, 8).Value = 2075.2856
$ws.Cells.Item(122, 9).Value = 1810.1111
$ws.Cells.Item(122, 11).Value = 5430.3333
$ws.Cells.Item(122, 13).Value = -2980.3333
$ws.Cells.Item(126, 8).Value = 2974.4285
$ws.Cells.Item(126, 9).Value = 2125.2856
$ws.Cells.Item(126, 10).Value = 3399
$ws.Cells.Item(126, 11).Value = 6375.8568
$ws.Cells.Item(126, 12).Value = 10197
$ws.Cells.Item(126, 13).Value = -3905.8568
$ws.Cells.Item(126, 14).Value = -15137

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(20, 8).Value = 31430.285
$ws.Cells.Item(20, 10).Value = 31430.285
$ws.Cells.Item(20, 12).Value = 31430.285
$ws.Cells.Item(20, 14).Value = -31882.285
$ws.Cells.Item(40, 8).Value = 3546.6775
$ws.Cells.Item(40, 9).Value = 2981.1667
$ws.Cells.Item(40, 11).Value = 2981.1667
$ws.Cells.Item(40, 13).Value = -2845.1667
$ws.Cells.Item(46, 8).Value = 5794.4116
$ws.Cells.Item(46, 9).Value = 2173
$ws.Cells.Item(46, 10).Value = 12433.667
$ws.Cells.Item(46, 11).Value = 2173
$ws.Cells.Item(46, 12).Value = 12433.667
$ws.Cells.Item(46, 13).Value = -1985
$ws.Cells.Item(46, 14).Value = -12809.667
$ws.Cells.Item(68, 8).Value = 5250
$ws.Cells.Item(68, 9).Value = 5250
$ws.Cells.Item(68, 11).Value = 5250
$ws.Cells.Item(68, 13).Value = -4501
$ws.Cells.Item(71, 8).Value = 5250
$ws.Cells.Item(71, 9).Value = 5250
$ws.Cells.Item(71, 11).Value = 26250
$ws.Cells.Item(71, 13).Value = -22506
$ws.Cells.Item(99, 8).Value = 0
$ws.Cells.Item(99, 10).Value = 0
$ws.Cells.Item(99, 12).Value = 0
$ws.Cells.Item(99, 14).ClearContents()
$ws.Cells.Item(132, 8).Value = 426901.22
$ws.Cells.Item(132, 9).Value = 11927.579
$ws.Cells.Item(132, 11).Value = 35782.737
$ws.Cells.Item(132, 13).Value = -33252.737

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(13, 8).Value = 2000
$ws.Cells.Item(13, 10).Value = 2000
$ws.Cells.Item(13, 12).Value = 2000
$ws.Cells.Item(13, 14).Value = -2280
$ws.Cells.Item(81, 8).Value = 12840
$ws.Cells.Item(81, 9).Value = 3067
$ws.Cells.Item(81, 11).Value = 6134
$ws.Cells.Item(81, 13).Value = -5073
$ws.Cells.Item(84, 8).Value = 12840
$ws.Cells.Item(84, 9).Value = 3067
$ws.Cells.Item(84, 11).Value = 30670
$ws.Cells.Item(84, 13).Value = -25366
$ws.Cells.Item(105, 8).Value = 30000
$ws.Cells.Item(105, 10).Value = 30000
$ws.Cells.Item(105, 12).Value = 30000
$ws.Cells.Item(105, 14).Value = -36988
